$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# SimTest09 - EEPROM functionality partially working: bug in Read routine.
# Update the "IST" score for that row from 0 to 2.5 ...
$ws.Range("D20").Value = 2.5

# ... and leave a comment explaining why.
$ws.Range("F20").Value = "TODO: Bug bei Read"

# Update the window view state (scroll position / selection) to reflect
# where the author was last working in the sheet.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C20").Select()
